$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Cotizacion N°"
$ws.Range("D1").Value = 500731
